$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at F; the old "Result" column (F) shifts to G.
$ws.Columns("F").Insert()

# Update header row: E1 becomes "Expected Output", new F1 is "Actual Output",
# G1 keeps "Result" (already shifted there by the insert).
$ws.Range("E1").Value = "Expected Output"
$ws.Range("F1").Value = "Actual Output"

# Populate the new "Actual Output" column (F) with the same text as the
# "Expected Output" column (E) for each data row, and fill in the two rows
# that previously had no Output value with "compiles" in both columns.
$ws.Range("E2").Value = "compiles"
$ws.Range("F2").Value = "compiles"

$ws.Range("F3").Value = $ws.Range("E3").Value()

$ws.Range("F4").Value = $ws.Range("E4").Value()

$ws.Range("F5").Value = $ws.Range("E5").Value()

$ws.Range("F6").Value = $ws.Range("E6").Value()

$ws.Range("E7").Value = "compiles"
$ws.Range("F7").Value = "compiles"

$ws.Range("F8").Value = $ws.Range("E8").Value()

$ws.Range("F9").Value = $ws.Range("E9").Value()

$ws.Range("F10").Value = $ws.Range("E10").Value()

$ws.Range("F11").Value = $ws.Range("E11").Value()

$ws.Range("F12").Value = $ws.Range("E12").Value()

# Widen the new "Actual Output" column to fit its header/content
# (target stored width ~27.285; closest attainable via this property).
$ws.Columns("F").ColumnWidth = 26.418

# Update the selection to match the post-edit state.
[void]$ws.Range("G14").Select()
